# paises.xlsx update - "Update countries & provincias Spain"
#
# The source feed refreshed the COVID country table:
#  - a handful of country rows (USA, Brasil, India, Canada, Ruanda/Libia
#    neighbourhood, Bermudas) got new totals;
#  - "Libia" and "Islas Malvinas" moved up in the underlying country
#    reference list, which (because the sheet's A column simply walks
#    that list row after row) shows up as the country label shifting on
#    a few rows while each row's own statistics (cols B:H) stay put;
#  - the "updated at" timestamp advanced from 01:12 to 02:29.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- "Datos actualizados..." footer (A1) -----------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Julio de 2020 a las 02:29"

# --- Country-list reshuffle, rows 131:134 (Libia now precedes Ruanda) ------
$ws.Range("A131").Value = "Libia"
$ws.Range("A132").Value = "Ruanda"
$ws.Range("A133").Value = "Jordania"
$ws.Range("A134").Value = "Letonia"

# --- Country-list reshuffle, rows 209:210 (Islas Malvinas now precedes Groenlandia) ---
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"

# --- Updated statistics -----------------------------------------------------
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 3096753
$ws.Range("C4").Value = 55111
$ws.Range("D4").Value = 1354717
$ws.Range("E4").Value = 1608065
$ws.Range("G4").Value = 992
$ws.Range("H4").Value = 133971

# Row 5 - Brasil
$ws.Range("B5").Value = 1674655
$ws.Range("C5").Value = 48584
$ws.Range("E5").Value = 535558
$ws.Range("G5").Value = 1312
$ws.Range("H5").Value = 66868

# Row 6 - India
$ws.Range("D6").Value = 457058
$ws.Range("E6").Value = 265770

# Row 23 - Canada
$ws.Range("B23").Value = 106167
$ws.Range("C23").Value = 232
$ws.Range("D23").Value = 69883
$ws.Range("E23").Value = 27573
$ws.Range("G23").Value = 18
$ws.Range("H23").Value = 8711

# Row 131 - now Libia (fresh numbers for the relocated country)
$ws.Range("B131").Value = 1182
$ws.Range("C131").Value = 65
$ws.Range("D131").Value = 295
$ws.Range("E131").Value = 852
$ws.Range("G131").Value = 1
$ws.Range("H131").Value = 35

# Row 132 - now Ruanda
$ws.Range("B132").Value = 1172
$ws.Range("C132").Value = 59
$ws.Range("D132").Value = 595
$ws.Range("E132").Value = 574
$ws.Range("H132").Value = 3

# Row 133 - now Jordania
$ws.Range("B133").Value = 1169
$ws.Range("C133").Value = 2
$ws.Range("D133").Value = 969
$ws.Range("E133").Value = 190
$ws.Range("H133").Value = 10

# Row 134 - now Letonia
$ws.Range("B134").Value = 1134
$ws.Range("C134").Value = 7
$ws.Range("D134").Value = 1008
$ws.Range("E134").Value = 96
$ws.Range("H134").Value = 30

# Row 176 - Bermudas
$ws.Range("B176").Value = 148
$ws.Range("C176").Value = 2
$ws.Range("E176").Value = 2
